$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four data rows (2-5) get rotated "up": the original row 2 values
# move down to row 5, while rows 3, 4, 5 shift up into 2, 3, 4.
# Columns D, M, N, O, P, S are the only ones that actually carry values
# that differ between rows (everything else is identical across rows
# 2-5 already), but we snapshot and rewrite them all the same way.

$cols = @("D", "M", "N", "O", "P", "S")

# Snapshot all original values for rows 2-5 before writing anything,
# so later writes don't clobber values we still need to read.
$orig = @{}
for ($r = 2; $r -le 5; $r++) {
    $orig[$r] = @{}
    foreach ($col in $cols) {
        $orig[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# New row order: row2<-old row3, row3<-old row4, row4<-old row5, row5<-old row2
$newOrder = @{ 2 = 3; 3 = 4; 4 = 5; 5 = 2 }

foreach ($r in $newOrder.Keys) {
    $src = $newOrder[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $orig[$src][$col]
    }
}
